$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: oxacycloheptadecan-2-one
$ws.Cells.Item(2, 1).Value = 'oxacycloheptadecan-2-one'
$ws.Cells.Item(2, 2).Value = 'oxacycloheptadecan-2-one'
$ws.Cells.Item(2, 3).Value = 'C16H30O2'
$ws.Cells.Item(2, 4).Value = 'C1CCCCCCCC(=O)OCCCCCCC1'
$ws.Cells.Item(2, 5).Value = 254.41
$ws.Cells.Item(2, 6).Value = 6.3
$ws.Cells.Item(2, 7).Value = 16
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = 30
$ws.Cells.Item(2, 10).Value = 2
$ws.Cells.Item(2, 11).Value = 0.7553791124562713
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = 0.1188632522306513
$ws.Cells.Item(2, 14).Value = 0.1257733579654888
$ws.Cells.Item(2, 15).Value = 14
$ws.Cells.Item(2, 16).Value = 0
$ws.Cells.Item(2, 17).Value = 0
$ws.Cells.Item(2, 18).Value = 0
$ws.Cells.Item(2, 19).Value = 0
$ws.Cells.Item(2, 20).Value = 1
$ws.Cells.Item(2, 21).Value = 0
$ws.Cells.Item(2, 22).Value = 0
$ws.Cells.Item(2, 23).Value = 0.771895758814512
$ws.Cells.Item(2, 24).Value = 0
$ws.Cells.Item(2, 25).Value = 0
$ws.Cells.Item(2, 26).Value = 0
$ws.Cells.Item(2, 27).Value = 0
$ws.Cells.Item(2, 28).Value = 0.2281199638378994
$ws.Cells.Item(2, 29).Value = 0
$ws.Cells.Item(2, 30).Value = 0
$ws.Cells.Item(2, 31).Value = 1.00018544935806

# Row 3: n-hexadecanoic acid
$ws.Cells.Item(3, 1).Value = 'n-hexadecanoic acid'
$ws.Cells.Item(3, 2).Value = 'hexadecanoic acid'
$ws.Cells.Item(3, 3).Value = 'C16H32O2'
$ws.Cells.Item(3, 4).Value = 'CCCCCCCCCCCCCCCC(=O)O'
$ws.Cells.Item(3, 5).Value = 256.42
$ws.Cells.Item(3, 6).Value = 6.4
$ws.Cells.Item(3, 7).Value = 16
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 9).Value = 32
$ws.Cells.Item(3, 10).Value = 2
$ws.Cells.Item(3, 11).Value = 0.7494579205990172
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 13).Value = 0.125793619842446
$ws.Cells.Item(3, 14).Value = 0.1247874580765931
$ws.Cells.Item(3, 15).Value = 15
$ws.Cells.Item(3, 16).Value = 0
$ws.Cells.Item(3, 17).Value = 0
$ws.Cells.Item(3, 18).Value = 0
$ws.Cells.Item(3, 19).Value = 1
$ws.Cells.Item(3, 20).Value = 0
$ws.Cells.Item(3, 21).Value = 0
$ws.Cells.Item(3, 22).Value = 0
$ws.Cells.Item(3, 23).Value = 0.8244793697839481
$ws.Cells.Item(3, 24).Value = 0
$ws.Cells.Item(3, 25).Value = 0
$ws.Cells.Item(3, 26).Value = 0
$ws.Cells.Item(3, 27).Value = 0.1755596287341081
$ws.Cells.Item(3, 28).Value = 0
$ws.Cells.Item(3, 29).Value = 0
$ws.Cells.Item(3, 30).Value = 0
$ws.Cells.Item(3, 31).Value = 1.00018544935806

# Row 4: tetradecanoic acid
$ws.Cells.Item(4, 1).Value = 'tetradecanoic acid'
$ws.Cells.Item(4, 2).Value = 'tetradecanoic acid'
$ws.Cells.Item(4, 3).Value = 'C14H28O2'
$ws.Cells.Item(4, 4).Value = 'CCCCCCCCCCCCCC(=O)O'
$ws.Cells.Item(4, 5).Value = 228.37
$ws.Cells.Item(4, 6).Value = 5.3
$ws.Cells.Item(4, 7).Value = 14
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 28
$ws.Cells.Item(4, 10).Value = 2
$ws.Cells.Item(4, 11).Value = 0.7363226343214958
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = 0.1235889127293427
$ws.Cells.Item(4, 14).Value = 0.1401147261023777
$ws.Cells.Item(4, 15).Value = 13
$ws.Cells.Item(4, 16).Value = 0
$ws.Cells.Item(4, 17).Value = 0
$ws.Cells.Item(4, 18).Value = 0
$ws.Cells.Item(4, 19).Value = 1
$ws.Cells.Item(4, 20).Value = 0
$ws.Cells.Item(4, 21).Value = 0
$ws.Cells.Item(4, 22).Value = 0
$ws.Cells.Item(4, 23).Value = 0.8029031834303979
$ws.Cells.Item(4, 24).Value = 0
$ws.Cells.Item(4, 25).Value = 0
$ws.Cells.Item(4, 26).Value = 0
$ws.Cells.Item(4, 27).Value = 0.1971230897228183
$ws.Cells.Item(4, 28).Value = 0
$ws.Cells.Item(4, 29).Value = 0
$ws.Cells.Item(4, 30).Value = 0
$ws.Cells.Item(4, 31).Value = 1.00018544935806

# Row 5: 2,4,5-trichlorophenol
$ws.Cells.Item(5, 1).Value = '2,4,5-trichlorophenol'
$ws.Cells.Item(5, 2).Value = '2,4,5-trichlorophenol'
$ws.Cells.Item(5, 3).Value = 'C6H3Cl3O'
$ws.Cells.Item(5, 4).Value = 'C1=C(C(=CC(=C1Cl)Cl)Cl)O'
$ws.Cells.Item(5, 5).Value = 197.4
$ws.Cells.Item(5, 6).Value = 3.7
$ws.Cells.Item(5, 7).Value = 6
$ws.Cells.Item(5, 8).Value = 3
$ws.Cells.Item(5, 9).Value = 3
$ws.Cells.Item(5, 10).Value = 1
$ws.Cells.Item(5, 11).Value = 0.3650759878419453
$ws.Cells.Item(5, 12).Value = 0.5387537993920973
$ws.Cells.Item(5, 13).Value = 0.01531914893617021
$ws.Cells.Item(5, 14).Value = 0.08104863221884498
$ws.Cells.Item(5, 15).Value = 0
$ws.Cells.Item(5, 16).Value = 6
$ws.Cells.Item(5, 17).Value = 3
$ws.Cells.Item(5, 18).Value = 1
$ws.Cells.Item(5, 19).Value = 0
$ws.Cells.Item(5, 20).Value = 0
$ws.Cells.Item(5, 21).Value = 3
$ws.Cells.Item(5, 22).Value = 0
$ws.Cells.Item(5, 23).Value = 0
$ws.Cells.Item(5, 24).Value = 0.3752887537993921
$ws.Cells.Item(5, 25).Value = 0.5387537993920973
$ws.Cells.Item(5, 26).Value = 0.0861550151975684
$ws.Cells.Item(5, 27).Value = 0
$ws.Cells.Item(5, 28).Value = 0
$ws.Cells.Item(5, 29).Value = 0.5387537993920973
$ws.Cells.Item(5, 30).Value = 0
$ws.Cells.Item(5, 31).Value = 1.00018544935806

# Row 6: oleic acid
$ws.Cells.Item(6, 1).Value = 'oleic acid'
$ws.Cells.Item(6, 2).Value = '(z)-octadec-9-enoic acid'
$ws.Cells.Item(6, 3).Value = 'C18H34O2'
$ws.Cells.Item(6, 4).Value = 'CCCCCCCCC=CCCCCCCCC(=O)O'
$ws.Cells.Item(6, 5).Value = 282.5
$ws.Cells.Item(6, 6).Value = 6.5
$ws.Cells.Item(6, 7).Value = 18
$ws.Cells.Item(6, 8).Value = 0
$ws.Cells.Item(6, 9).Value = 34
$ws.Cells.Item(6, 10).Value = 2
$ws.Cells.Item(6, 11).Value = 0.7653026548672566
$ws.Cells.Item(6, 12).Value = 0
$ws.Cells.Item(6, 13).Value = 0.121316814159292
$ws.Cells.Item(6, 14).Value = 0.1132672566371681
$ws.Cells.Item(6, 15).Value = 17
$ws.Cells.Item(6, 16).Value = 0
$ws.Cells.Item(6, 17).Value = 0
$ws.Cells.Item(6, 18).Value = 0
$ws.Cells.Item(6, 19).Value = 1
$ws.Cells.Item(6, 20).Value = 0
$ws.Cells.Item(6, 21).Value = 0
$ws.Cells.Item(6, 22).Value = 0
$ws.Cells.Item(6, 23).Value = 0.8405345132743363
$ws.Cells.Item(6, 24).Value = 0
$ws.Cells.Item(6, 25).Value = 0
$ws.Cells.Item(6, 26).Value = 0
$ws.Cells.Item(6, 27).Value = 0.1593522123893805
$ws.Cells.Item(6, 28).Value = 0
$ws.Cells.Item(6, 29).Value = 0
$ws.Cells.Item(6, 30).Value = 0
$ws.Cells.Item(6, 31).Value = 1.00018544935806

# Row 7: 2-butanone
$ws.Cells.Item(7, 1).Value = '2-butanone'
$ws.Cells.Item(7, 2).Value = 'butan-2-one'
$ws.Cells.Item(7, 3).Value = 'C4H8O'
$ws.Cells.Item(7, 4).Value = 'CCC(=O)C'
$ws.Cells.Item(7, 5).Value = 72.11
$ws.Cells.Item(7, 6).Value = 0.3
$ws.Cells.Item(7, 7).Value = 4
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(7, 9).Value = 8
$ws.Cells.Item(7, 10).Value = 1
$ws.Cells.Item(7, 11).Value = 0.6662598807377618
$ws.Cells.Item(7, 12).Value = 0
$ws.Cells.Item(7, 13).Value = 0.1118291499098599
$ws.Cells.Item(7, 14).Value = 0.2218693662460131
$ws.Cells.Item(7, 15).Value = 1
$ws.Cells.Item(7, 16).Value = 0
$ws.Cells.Item(7, 17).Value = 0
$ws.Cells.Item(7, 18).Value = 0
$ws.Cells.Item(7, 19).Value = 0
$ws.Cells.Item(7, 20).Value = 0
$ws.Cells.Item(7, 21).Value = 0
$ws.Cells.Item(7, 22).Value = 1
$ws.Cells.Item(7, 23).Value = 0.2085009014006379
$ws.Cells.Item(7, 24).Value = 0
$ws.Cells.Item(7, 25).Value = 0
$ws.Cells.Item(7, 26).Value = 0
$ws.Cells.Item(7, 27).Value = 0
$ws.Cells.Item(7, 28).Value = 0
$ws.Cells.Item(7, 29).Value = 0
$ws.Cells.Item(7, 30).Value = 0.7914574954929968
$ws.Cells.Item(7, 31).Value = 1.00018544935806

# Row 8: 2-cyclopenten-1-one, 2-methyl-
$ws.Cells.Item(8, 1).Value = '2-cyclopenten-1-one, 2-methyl-'
$ws.Cells.Item(8, 2).Value = '2-methylcyclopent-2-en-1-one'
$ws.Cells.Item(8, 3).Value = 'C6H8O'
$ws.Cells.Item(8, 4).Value = 'CC1=CCCC1=O'
$ws.Cells.Item(8, 5).Value = 96.13
$ws.Cells.Item(8, 6).Value = 0.9
$ws.Cells.Item(8, 7).Value = 6
$ws.Cells.Item(8, 8).Value = 0
$ws.Cells.Item(8, 9).Value = 8
$ws.Cells.Item(8, 10).Value = 1
$ws.Cells.Item(8, 11).Value = 0.7496723187350464
$ws.Cells.Item(8, 12).Value = 0
$ws.Cells.Item(8, 13).Value = 0.08388640382814938
$ws.Cells.Item(8, 14).Value = 0.1664308748569645
$ws.Cells.Item(8, 15).Value = 3
$ws.Cells.Item(8, 16).Value = 0
$ws.Cells.Item(8, 17).Value = 0
$ws.Cells.Item(8, 18).Value = 0
$ws.Cells.Item(8, 19).Value = 0
$ws.Cells.Item(8, 20).Value = 0
$ws.Cells.Item(8, 21).Value = 0
$ws.Cells.Item(8, 22).Value = 1
$ws.Cells.Item(8, 23).Value = 0.4377509622386352
$ws.Cells.Item(8, 24).Value = 0
$ws.Cells.Item(8, 25).Value = 0
$ws.Cells.Item(8, 26).Value = 0
$ws.Cells.Item(8, 27).Value = 0
$ws.Cells.Item(8, 28).Value = 0
$ws.Cells.Item(8, 29).Value = 0
$ws.Cells.Item(8, 30).Value = 0.5936960366170811
$ws.Cells.Item(8, 31).Value = 1.00018544935806

# Row 9: 2-methylcyclopent-2-en-1-one
$ws.Cells.Item(9, 1).Value = '2-methylcyclopent-2-en-1-one'
$ws.Cells.Item(9, 2).Value = '2-methylcyclopent-2-en-1-one'
$ws.Cells.Item(9, 3).Value = 'C6H8O'
$ws.Cells.Item(9, 4).Value = 'CC1=CCCC1=O'
$ws.Cells.Item(9, 5).Value = 96.13
$ws.Cells.Item(9, 6).Value = 0.9
$ws.Cells.Item(9, 7).Value = 6
$ws.Cells.Item(9, 8).Value = 0
$ws.Cells.Item(9, 9).Value = 8
$ws.Cells.Item(9, 10).Value = 1
$ws.Cells.Item(9, 11).Value = 0.7496723187350464
$ws.Cells.Item(9, 12).Value = 0
$ws.Cells.Item(9, 13).Value = 0.08388640382814938
$ws.Cells.Item(9, 14).Value = 0.1664308748569645
$ws.Cells.Item(9, 15).Value = 3
$ws.Cells.Item(9, 16).Value = 0
$ws.Cells.Item(9, 17).Value = 0
$ws.Cells.Item(9, 18).Value = 0
$ws.Cells.Item(9, 19).Value = 0
$ws.Cells.Item(9, 20).Value = 0
$ws.Cells.Item(9, 21).Value = 0
$ws.Cells.Item(9, 22).Value = 1
$ws.Cells.Item(9, 23).Value = 0.4377509622386352
$ws.Cells.Item(9, 24).Value = 0
$ws.Cells.Item(9, 25).Value = 0
$ws.Cells.Item(9, 26).Value = 0
$ws.Cells.Item(9, 27).Value = 0
$ws.Cells.Item(9, 28).Value = 0
$ws.Cells.Item(9, 29).Value = 0
$ws.Cells.Item(9, 30).Value = 0.5936960366170811
$ws.Cells.Item(9, 31).Value = 1.00018544935806

# Row 10: (z)-octadec-9-enoic acid
$ws.Cells.Item(10, 1).Value = '(z)-octadec-9-enoic acid'
$ws.Cells.Item(10, 2).Value = '(z)-octadec-9-enoic acid'
$ws.Cells.Item(10, 3).Value = 'C18H34O2'
$ws.Cells.Item(10, 4).Value = 'CCCCCCCCC=CCCCCCCCC(=O)O'
$ws.Cells.Item(10, 5).Value = 282.5
$ws.Cells.Item(10, 6).Value = 6.5
$ws.Cells.Item(10, 7).Value = 18
$ws.Cells.Item(10, 8).Value = 0
$ws.Cells.Item(10, 9).Value = 34
$ws.Cells.Item(10, 10).Value = 2
$ws.Cells.Item(10, 11).Value = 0.7653026548672566
$ws.Cells.Item(10, 12).Value = 0
$ws.Cells.Item(10, 13).Value = 0.121316814159292
$ws.Cells.Item(10, 14).Value = 0.1132672566371681
$ws.Cells.Item(10, 15).Value = 17
$ws.Cells.Item(10, 16).Value = 0
$ws.Cells.Item(10, 17).Value = 0
$ws.Cells.Item(10, 18).Value = 0
$ws.Cells.Item(10, 19).Value = 1
$ws.Cells.Item(10, 20).Value = 0
$ws.Cells.Item(10, 21).Value = 0
$ws.Cells.Item(10, 22).Value = 0
$ws.Cells.Item(10, 23).Value = 0.8405345132743363
$ws.Cells.Item(10, 24).Value = 0
$ws.Cells.Item(10, 25).Value = 0
$ws.Cells.Item(10, 26).Value = 0
$ws.Cells.Item(10, 27).Value = 0.1593522123893805
$ws.Cells.Item(10, 28).Value = 0
$ws.Cells.Item(10, 29).Value = 0
$ws.Cells.Item(10, 30).Value = 0
$ws.Cells.Item(10, 31).Value = 1.00018544935806

# Row 11: 2,5-hexanedione
$ws.Cells.Item(11, 1).Value = '2,5-hexanedione'
$ws.Cells.Item(11, 2).Value = 'hexane-2,5-dione'
$ws.Cells.Item(11, 3).Value = 'C6H10O2'
$ws.Cells.Item(11, 4).Value = 'CC(=O)CCC(=O)C'
$ws.Cells.Item(11, 5).Value = 114.14
$ws.Cells.Item(11, 6).Value = -0.3
$ws.Cells.Item(11, 7).Value = 6
$ws.Cells.Item(11, 8).Value = 0
$ws.Cells.Item(11, 9).Value = 10
$ws.Cells.Item(11, 10).Value = 2
$ws.Cells.Item(11, 11).Value = 0.6313825127036973
$ws.Cells.Item(11, 12).Value = 0
$ws.Cells.Item(11, 13).Value = 0.08831259856316805
$ws.Cells.Item(11, 14).Value = 0.2803399334151043
$ws.Cells.Item(11, 15).Value = 0
$ws.Cells.Item(11, 16).Value = 0
$ws.Cells.Item(11, 17).Value = 0
$ws.Cells.Item(11, 18).Value = 0
$ws.Cells.Item(11, 19).Value = 0
$ws.Cells.Item(11, 20).Value = 0
$ws.Cells.Item(11, 21).Value = 0
$ws.Cells.Item(11, 22).Value = 2
$ws.Cells.Item(11, 23).Value = 0
$ws.Cells.Item(11, 24).Value = 0
$ws.Cells.Item(11, 25).Value = 0
$ws.Cells.Item(11, 26).Value = 0
$ws.Cells.Item(11, 27).Value = 0
$ws.Cells.Item(11, 28).Value = 0
$ws.Cells.Item(11, 29).Value = 0
$ws.Cells.Item(11, 30).Value = 1.000035044681969
$ws.Cells.Item(11, 31).Value = 1.00018544935806

# Row 12: 1-hexene, 4,5-dimethyl-
$ws.Cells.Item(12, 1).Value = '1-hexene, 4,5-dimethyl-'
$ws.Cells.Item(12, 2).Value = '4,5-dimethylhex-1-ene'
$ws.Cells.Item(12, 3).Value = 'C8H16'
$ws.Cells.Item(12, 4).Value = 'CC(C)C(C)CC=C'
$ws.Cells.Item(12, 5).Value = 112.21
$ws.Cells.Item(12, 6).Value = 3.5
$ws.Cells.Item(12, 7).Value = 8
$ws.Cells.Item(12, 8).Value = 0
$ws.Cells.Item(12, 9).Value = 16
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 11).Value = 0.8563229658675697
$ws.Cells.Item(12, 12).Value = 0
$ws.Cells.Item(12, 13).Value = 0.1437305053025577
$ws.Cells.Item(12, 14).Value = 0
$ws.Cells.Item(12, 15).Value = 8
$ws.Cells.Item(12, 16).Value = 0
$ws.Cells.Item(12, 17).Value = 0
$ws.Cells.Item(12, 18).Value = 0
$ws.Cells.Item(12, 19).Value = 0
$ws.Cells.Item(12, 20).Value = 0
$ws.Cells.Item(12, 21).Value = 0
$ws.Cells.Item(12, 22).Value = 0
$ws.Cells.Item(12, 23).Value = 1.000053471170127
$ws.Cells.Item(12, 24).Value = 0
$ws.Cells.Item(12, 25).Value = 0
$ws.Cells.Item(12, 26).Value = 0
$ws.Cells.Item(12, 27).Value = 0
$ws.Cells.Item(12, 28).Value = 0
$ws.Cells.Item(12, 29).Value = 0
$ws.Cells.Item(12, 30).Value = 0
$ws.Cells.Item(12, 31).Value = 1.00018544935806

# Row 13: trans-2-pentenoic acid
$ws.Cells.Item(13, 1).Value = 'trans-2-pentenoic acid'
$ws.Cells.Item(13, 2).Value = '(e)-pent-2-enoic acid'
$ws.Cells.Item(13, 3).Value = 'C5H8O2'
$ws.Cells.Item(13, 4).Value = 'CCC=CC(=O)O'
$ws.Cells.Item(13, 5).Value = 100.12
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 5
$ws.Cells.Item(13, 8).Value = 0
$ws.Cells.Item(13, 9).Value = 8
$ws.Cells.Item(13, 10).Value = 2
$ws.Cells.Item(13, 11).Value = 0.5998302037554933
$ws.Cells.Item(13, 12).Value = 0
$ws.Cells.Item(13, 13).Value = 0.08054334798242109
$ws.Cells.Item(13, 14).Value = 0.3195964842189373
$ws.Cells.Item(13, 15).Value = 4
$ws.Cells.Item(13, 16).Value = 0
$ws.Cells.Item(13, 17).Value = 0
$ws.Cells.Item(13, 18).Value = 0
$ws.Cells.Item(13, 19).Value = 1
$ws.Cells.Item(13, 20).Value = 0
$ws.Cells.Item(13, 21).Value = 0
$ws.Cells.Item(13, 22).Value = 0
$ws.Cells.Item(13, 23).Value = 0.5503395924890131
$ws.Cells.Item(13, 24).Value = 0
$ws.Cells.Item(13, 25).Value = 0
$ws.Cells.Item(13, 26).Value = 0
$ws.Cells.Item(13, 27).Value = 0.4496304434678386
$ws.Cells.Item(13, 28).Value = 0
$ws.Cells.Item(13, 29).Value = 0
$ws.Cells.Item(13, 30).Value = 0
$ws.Cells.Item(13, 31).Value = 1.00018544935806

# Row 14: n-decanoic acid
$ws.Cells.Item(14, 1).Value = 'n-decanoic acid'
$ws.Cells.Item(14, 2).Value = 'decanoic acid'
$ws.Cells.Item(14, 3).Value = 'C10H20O2'
$ws.Cells.Item(14, 4).Value = 'CCCCCCCCCC(=O)O'
$ws.Cells.Item(14, 5).Value = 172.26
$ws.Cells.Item(14, 6).Value = 4.1
$ws.Cells.Item(14, 7).Value = 10
$ws.Cells.Item(14, 8).Value = 0
$ws.Cells.Item(14, 9).Value = 20
$ws.Cells.Item(14, 10).Value = 2
$ws.Cells.Item(14, 11).Value = 0.6972599558806455
$ws.Cells.Item(14, 12).Value = 0
$ws.Cells.Item(14, 13).Value = 0.1170323928944619
$ws.Cells.Item(14, 14).Value = 0.1857540926506444
$ws.Cells.Item(14, 15).Value = 9
$ws.Cells.Item(14, 16).Value = 0
$ws.Cells.Item(14, 17).Value = 0
$ws.Cells.Item(14, 18).Value = 0
$ws.Cells.Item(14, 19).Value = 1
$ws.Cells.Item(14, 20).Value = 0
$ws.Cells.Item(14, 21).Value = 0
$ws.Cells.Item(14, 22).Value = 0
$ws.Cells.Item(14, 23).Value = 0.7387147335423198
$ws.Cells.Item(14, 24).Value = 0
$ws.Cells.Item(14, 25).Value = 0
$ws.Cells.Item(14, 26).Value = 0
$ws.Cells.Item(14, 27).Value = 0.2613317078834321
$ws.Cells.Item(14, 28).Value = 0
$ws.Cells.Item(14, 29).Value = 0
$ws.Cells.Item(14, 30).Value = 0
$ws.Cells.Item(14, 31).Value = 1.00018544935806

# Row 15: hexadecanoic acid
$ws.Cells.Item(15, 1).Value = 'hexadecanoic acid'
$ws.Cells.Item(15, 2).Value = 'hexadecanoic acid'
$ws.Cells.Item(15, 3).Value = 'C16H32O2'
$ws.Cells.Item(15, 4).Value = 'CCCCCCCCCCCCCCCC(=O)O'
$ws.Cells.Item(15, 5).Value = 256.42
$ws.Cells.Item(15, 6).Value = 6.4
$ws.Cells.Item(15, 7).Value = 16
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 9).Value = 32
$ws.Cells.Item(15, 10).Value = 2
$ws.Cells.Item(15, 11).Value = 0.7494579205990172
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 13).Value = 0.125793619842446
$ws.Cells.Item(15, 14).Value = 0.1247874580765931
$ws.Cells.Item(15, 15).Value = 15
$ws.Cells.Item(15, 16).Value = 0
$ws.Cells.Item(15, 17).Value = 0
$ws.Cells.Item(15, 18).Value = 0
$ws.Cells.Item(15, 19).Value = 1
$ws.Cells.Item(15, 20).Value = 0
$ws.Cells.Item(15, 21).Value = 0
$ws.Cells.Item(15, 22).Value = 0
$ws.Cells.Item(15, 23).Value = 0.8244793697839481
$ws.Cells.Item(15, 24).Value = 0
$ws.Cells.Item(15, 25).Value = 0
$ws.Cells.Item(15, 26).Value = 0
$ws.Cells.Item(15, 27).Value = 0.1755596287341081
$ws.Cells.Item(15, 28).Value = 0
$ws.Cells.Item(15, 29).Value = 0
$ws.Cells.Item(15, 30).Value = 0
$ws.Cells.Item(15, 31).Value = 1.00018544935806

# Row 16: phenol
$ws.Cells.Item(16, 1).Value = 'phenol'
$ws.Cells.Item(16, 2).Value = 'phenol'
$ws.Cells.Item(16, 3).Value = 'C6H6O'
$ws.Cells.Item(16, 4).Value = 'C1=CC=C(C=C1)O'
$ws.Cells.Item(16, 5).Value = 94.11
$ws.Cells.Item(16, 6).Value = 1.5
$ws.Cells.Item(16, 7).Value = 6
$ws.Cells.Item(16, 8).Value = 0
$ws.Cells.Item(16, 9).Value = 6
$ws.Cells.Item(16, 10).Value = 1
$ws.Cells.Item(16, 11).Value = 0.765763468281798
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 13).Value = 0.06426522154925088
$ws.Cells.Item(16, 14).Value = 0.1700031877590054
$ws.Cells.Item(16, 15).Value = 0
$ws.Cells.Item(16, 16).Value = 6
$ws.Cells.Item(16, 17).Value = 0
$ws.Cells.Item(16, 18).Value = 1
$ws.Cells.Item(16, 19).Value = 0
$ws.Cells.Item(16, 20).Value = 0
$ws.Cells.Item(16, 21).Value = 0
$ws.Cells.Item(16, 22).Value = 0
$ws.Cells.Item(16, 23).Value = 0
$ws.Cells.Item(16, 24).Value = 0.8193178195728402
$ws.Cells.Item(16, 25).Value = 0
$ws.Cells.Item(16, 26).Value = 0.1807140580172139
$ws.Cells.Item(16, 27).Value = 0
$ws.Cells.Item(16, 28).Value = 0
$ws.Cells.Item(16, 29).Value = 0
$ws.Cells.Item(16, 30).Value = 0
$ws.Cells.Item(16, 31).Value = 1.00018544935806

# Row 17: 9,12-octadecadienoic acid (z,z)-
$ws.Cells.Item(17, 1).Value = '9,12-octadecadienoic acid (z,z)-'
$ws.Cells.Item(17, 2).Value = '(9z,12z)-octadeca-9,12-dienoic acid'
$ws.Cells.Item(17, 3).Value = 'C18H32O2'
$ws.Cells.Item(17, 4).Value = 'CCCCCC=CCC=CCCCCCCCC(=O)O'
$ws.Cells.Item(17, 5).Value = 280.4
$ws.Cells.Item(17, 6).Value = 6.8
$ws.Cells.Item(17, 7).Value = 18
$ws.Cells.Item(17, 8).Value = 0
$ws.Cells.Item(17, 9).Value = 32
$ws.Cells.Item(17, 10).Value = 2
$ws.Cells.Item(17, 11).Value = 0.7710342368045648
$ws.Cells.Item(17, 12).Value = 0
$ws.Cells.Item(17, 13).Value = 0.1150356633380885
$ws.Cells.Item(17, 14).Value = 0.1141155492154066
$ws.Cells.Item(17, 15).Value = 17
$ws.Cells.Item(17, 16).Value = 0
$ws.Cells.Item(17, 17).Value = 0
$ws.Cells.Item(17, 18).Value = 0
$ws.Cells.Item(17, 19).Value = 1
$ws.Cells.Item(17, 20).Value = 0
$ws.Cells.Item(17, 21).Value = 0
$ws.Cells.Item(17, 22).Value = 0
$ws.Cells.Item(17, 23).Value = 0.8396398002853066
$ws.Cells.Item(17, 24).Value = 0
$ws.Cells.Item(17, 25).Value = 0
$ws.Cells.Item(17, 26).Value = 0
$ws.Cells.Item(17, 27).Value = 0.1605456490727532
$ws.Cells.Item(17, 28).Value = 0
$ws.Cells.Item(17, 29).Value = 0
$ws.Cells.Item(17, 30).Value = 0
$ws.Cells.Item(17, 31).Value = 1.00018544935806

# Row 18: (9z,12z)-octadeca-9,12-dienoic acid
$ws.Cells.Item(18, 1).Value = '(9z,12z)-octadeca-9,12-dienoic acid'
$ws.Cells.Item(18, 2).Value = '(9z,12z)-octadeca-9,12-dienoic acid'
$ws.Cells.Item(18, 3).Value = 'C18H32O2'
$ws.Cells.Item(18, 4).Value = 'CCCCCC=CCC=CCCCCCCCC(=O)O'
$ws.Cells.Item(18, 5).Value = 280.4
$ws.Cells.Item(18, 6).Value = 6.8
$ws.Cells.Item(18, 7).Value = 18
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 9).Value = 32
$ws.Cells.Item(18, 10).Value = 2
$ws.Cells.Item(18, 11).Value = 0.7710342368045648
$ws.Cells.Item(18, 12).Value = 0
$ws.Cells.Item(18, 13).Value = 0.1150356633380885
$ws.Cells.Item(18, 14).Value = 0.1141155492154066
$ws.Cells.Item(18, 15).Value = 17
$ws.Cells.Item(18, 16).Value = 0
$ws.Cells.Item(18, 17).Value = 0
$ws.Cells.Item(18, 18).Value = 0
$ws.Cells.Item(18, 19).Value = 1
$ws.Cells.Item(18, 20).Value = 0
$ws.Cells.Item(18, 21).Value = 0
$ws.Cells.Item(18, 22).Value = 0
$ws.Cells.Item(18, 23).Value = 0.8396398002853066
$ws.Cells.Item(18, 24).Value = 0
$ws.Cells.Item(18, 25).Value = 0
$ws.Cells.Item(18, 26).Value = 0
$ws.Cells.Item(18, 27).Value = 0.1605456490727532
$ws.Cells.Item(18, 28).Value = 0
$ws.Cells.Item(18, 29).Value = 0
$ws.Cells.Item(18, 30).Value = 0
$ws.Cells.Item(18, 31).Value = 1.00018544935806
